$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.186.94'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.643.57'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.61'
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.19'
$ws.Range('E6').Value = '  +0.93%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.144'
$ws.Range('E9').Value = '  +5.84%  '
$ws.Range('E10').Value = '  -0.76%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.352'
$ws.Range('E12').Value = '  +1.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.09'
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000191'
$ws.Range('E14').Value = '  +1.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.127.55'
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '68.121.45'
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('E17').Value = '  +0.61%  '
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '363.17'
$ws.Range('E19').Value = '  -1.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.44'
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('E21').Value = '  +3.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.83'
$ws.Range('E22').Value = '  -0.78%  '
$ws.Range('E23').Value = '  -1.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.02'
$ws.Range('E24').Value = '  +2.50%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('E26').Value = '  -1.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.790.31'
$ws.Range('E27').Value = '  +0.95%  '
$ws.Range('E28').Value = '  +0.66%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '556.24'
$ws.Range('E30').Value = '  -3.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.01'
$ws.Range('E31').Value = '  +0.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.41'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.85'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('E34').Value = '  +1.04%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.56'
$ws.Range('E36').Value = '  +2.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '160.99'
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('E38').Value = '  +0.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.373'
$ws.Range('E39').Value = '  +1.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.87'
$ws.Range('E40').Value = '  -2.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.33'
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0₆0340'
$ws.Range('E42').Value = '  +4.99%  '
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.77'
$ws.Range('E43').Value = '  +0.82%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.63'
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.36'
$ws.Range('E46').Value = '  -0.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '158.27'
$ws.Range('E47').Value = '  +1.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.73'
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '21.98'
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0786'
$ws.Range('E51').Value = '  +0.74%  '
